$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue 2 4 "27.359.06"
Set-TextValue 2 5 "  +1.26%  "

# Row 3
Set-TextValue 3 4 "1.859.63"
Set-TextValue 3 5 "  +2.04%  "

# Row 4
Set-TextValue 4 4 "1.005"
Set-TextValue 4 5 "  -0.28%  "

# Row 5
Set-TextValue 5 4 "314.44"
Set-TextValue 5 5 "  +1.60%  "

# Row 6
Set-TextValue 6 4 "1.004"
Set-TextValue 6 5 "  -0.30%  "

# Row 7
Set-TextValue 7 4 "0.4642"
Set-TextValue 7 5 "  +0.08%  "

# Row 8
Set-TextValue 8 4 "0.3718"
Set-TextValue 8 5 "  +1.56%  "

# Row 9
Set-TextValue 9 4 "0.07343"
Set-TextValue 9 5 "  +1.47%  "

# Row 10
Set-TextValue 10 4 "0.8829"
Set-TextValue 10 5 "  +2.76%  "

# Row 11
Set-TextValue 11 4 "0.07908"
Set-TextValue 11 5 "  +3.48%  "

# Row 12
Set-TextValue 12 4 "19.86"
Set-TextValue 12 5 "  +0.04%  "

# Row 13
Set-TextValue 13 4 "1.943.21"
Set-TextValue 13 5 "  +5.53%  "

# Row 14
Set-TextValue 14 4 "5.400"
Set-TextValue 14 5 "  +1.37%  "

# Row 15
Set-TextValue 15 4 "6.579"
Set-TextValue 15 5 "  +0.93%  "

# Row 16
Set-TextValue 16 4 "92.14"
Set-TextValue 16 5 "  +0.49%  "

# Row 17
Set-TextValue 17 4 "1.005"
Set-TextValue 17 5 "  -0.23%  "

# Row 18
Set-TextValue 18 4 "0.000008870"
Set-TextValue 18 5 "  +2.65%  "

# Row 19
Set-TextValue 19 5 "  -0.44%  "

# Row 20
Set-TextValue 20 4 "14.83"
Set-TextValue 20 5 "  +2.46%  "

# Row 21
Set-TextValue 21 4 "27.400.57"
Set-TextValue 21 5 "  +0.73%  "

# Row 22
Set-TextValue 22 4 "5.137"
Set-TextValue 22 5 "  -0.21%  "

# Row 23
Set-TextValue 23 5 "  +0.13%  "

# Row 24
Set-TextValue 24 4 "2.154.99"
Set-TextValue 24 5 "  +8.88%  "

# Row 25
Set-TextValue 25 4 "1.898"
Set-TextValue 25 5 "  +3.12%  "

# Row 26
Set-TextValue 26 4 "153.05"
Set-TextValue 26 5 "  +0.85%  "

# Row 27
Set-TextValue 27 4 "18.47"
Set-TextValue 27 5 "  +1.71%  "

# Row 28
Set-TextValue 28 4 "2.073"
Set-TextValue 28 5 "  +1.45%  "

# Row 29
Set-TextValue 29 4 "5.127"
Set-TextValue 29 5 "  +0.50%  "

# Row 30
Set-TextValue 30 4 "116.37"
Set-TextValue 30 5 "  +0.91%  "

# Row 31
Set-TextValue 31 4 "0.08893"
Set-TextValue 31 5 "  +0.62%  "

# Row 32
Set-TextValue 32 2 "HuobiToken"
Set-TextValue 32 3 "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue 32 4 "3.026"
Set-TextValue 32 5 "  +2.39%  "

# Row 33
Set-TextValue 33 2 "ImmutableX"
Set-TextValue 33 3 "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue 33 4 "0.7557"
Set-TextValue 33 5 "  +5.05%  "

# Row 34
Set-TextValue 34 4 "1.160"
Set-TextValue 34 5 "  +2.78%  "

# Row 35
Set-TextValue 35 4 "4.483"
Set-TextValue 35 5 "  +1.26%  "

# Row 36
Set-TextValue 36 4 "2.642"
Set-TextValue 36 5 "  +9.85%  "

# Row 37
Set-TextValue 37 4 "0.01963"
Set-TextValue 37 5 "  +1.89%  "

# Row 38
Set-TextValue 38 4 "1.076"
Set-TextValue 38 5 "  -0.06%  "

# Row 39
Set-TextValue 39 2 "MXToken"
Set-TextValue 39 3 "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue 39 4 "2.980"
Set-TextValue 39 5 "  +1.79%  "

# Row 40
Set-TextValue 40 2 "Hedera"
Set-TextValue 40 3 "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue 40 4 "0.05237"
Set-TextValue 40 5 "  -0.27%  "

# Row 41
Set-TextValue 41 4 "7.085"
Set-TextValue 41 5 "  -0.69%  "

# Row 42
Set-TextValue 42 4 "0.5163"
Set-TextValue 42 5 "  -0.13%  "

# Row 43
Set-TextValue 43 4 "0.1644"
Set-TextValue 43 5 "  +1.06%  "

# Row 44
Set-TextValue 44 4 "8.332"
Set-TextValue 44 5 "  +2.16%  "

# Row 45
Set-TextValue 45 4 "0.4836"
Set-TextValue 45 5 "  +0.76%  "

# Row 46
Set-TextValue 46 4 "10.32"
Set-TextValue 46 5 "  +1.53%  "

# Row 47
Set-TextValue 47 4 "1.004"
Set-TextValue 47 5 "  -0.38%  "

# Row 48
Set-TextValue 48 4 "103.33"
Set-TextValue 48 5 "  +0.62%  "

# Row 49
Set-TextValue 49 4 "1.651"
Set-TextValue 49 5 "  +2.10%  "

# Row 50
Set-TextValue 50 4 "0.06236"
Set-TextValue 50 5 "  -0.30%  "

# Row 51
Set-TextValue 51 4 "65.55"
Set-TextValue 51 5 "  +1.86%  "
